$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.757.28'
$ws.Range("E2").Value = '  -2.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.482.41'
$ws.Range("E3").Value = '  -2.91%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.74'
$ws.Range("E5").Value = '  -3.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.33'
$ws.Range("E6").Value = '  -4.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.487.27'
$ws.Range("E7").Value = '  -2.73%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.124'
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.07'
$ws.Range("E11").Value = '  -1.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.385'
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.059.81'
$ws.Range("E13").Value = '  -3.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.74'
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000178'
$ws.Range("E15").Value = '  -4.66%  '
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.473.03'
$ws.Range("E17").Value = '  -3.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.771.10'
$ws.Range("E18").Value = '  -3.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.94'
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.30'
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("E21").Value = '  -4.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.53'
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.576'
$ws.Range("E23").Value = '  -2.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.614.31'
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '72.70'
$ws.Range("E26").Value = '  -2.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000109'
$ws.Range("E27").Value = '  -7.86%  '
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.42'
$ws.Range("E30").Value = '  -8.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.24'
$ws.Range("E31").Value = '  -6.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.16'
$ws.Range("E32").Value = '  -5.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.478.71'
$ws.Range("E33").Value = '  -3.09%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.70'
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("E36").Value = '  -3.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.30'
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.94'
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.56'
$ws.Range("E39").Value = '  -3.33%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '168.87'
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0805'
$ws.Range("E41").Value = '  -3.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.809'
$ws.Range("E42").Value = '  -3.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.05'
$ws.Range("E43").Value = '  -3.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.76'
$ws.Range("E45").Value = '  -3.17%  '
$ws.Range("E46").Value = '  -6.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.34'
$ws.Range("E47").Value = '  -4.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.63'
$ws.Range("E48").Value = '  -4.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.87'
$ws.Range("E49").Value = '  -2.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.421.15'
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("E51").Value = '  -1.83%  '
